# Rename the three header/footer logo pictures' internal image names:
#   - footer (primary)    id=2 : image2.png -> image1.png
#   - footer (first page) id=3 : image2.png -> image1.png
#   - header (first page) id=1 : image1.jpg -> image2.jpg
#
# InlineShape objects don't reliably carry a live handle for property
# writes straight off a Range/HeaderFooter collection in this host, so each
# picture is selected first and then re-acquired from $word.Selection
# (the pattern the host's own stale-handle hint points at) before the
# .Name assignment, which is what actually lands the rename into the
# picture's docPr/name.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-FooterLogo($footerIndex, $newName) {
    $footer = $sec.Footers.Item($footerIndex)
    $pic = $footer.Range.InlineShapes.Item(1)
    $pic.Select()
    $selected = $word.Selection.InlineShapes.Item(1)
    $selected.Name = $newName
}

function Rename-HeaderLogo($headerIndex, $newName) {
    $header = $sec.Headers.Item($headerIndex)
    $pic = $header.Range.InlineShapes.Item(1)
    $pic.Select()
    $selected = $word.Selection.InlineShapes.Item(1)
    $selected.Name = $newName
}

# Primary footer (wdHeaderFooterPrimary = 1) -> Pearson logo id="2"
Rename-FooterLogo 1 "image1.png"

# First-page footer (wdHeaderFooterFirstPage = 2) -> Pearson logo id="3"
Rename-FooterLogo 2 "image1.png"

# First-page header (wdHeaderFooterFirstPage = 2) -> BTEC logo id="1"
Rename-HeaderLogo 2 "image2.jpg"

Write-Host "Renamed footer/header logo pictures."
